# TestcaseSuite.xlsx — "First 3 working cases" 17/6
#
# TC_01 (row 2): append an "assert" step (keyword/data/objectName columns
# become comma-separated lists) to check the landing page URL after login.
#
# TC_02 (row 3, new): verify clicking "Create" navigates to the
# "Create New Client" screen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update TC_01 (row 2) ---
$ws.Range("C2").Value = "openbrowser,assert"
$ws.Range("D2").Value = "https://adva-pro-dev01.paradigmcentral.com,https://adva-pro-dev01.paradigmcentral.com/#/providers/prospects/list"
$ws.Range("E2").Value = "no value,url"
$ws.Range("F2").Value = "yes"

# --- Add TC_02 (row 3) ---
$ws.Range("A3").Value = "TC_02"
$ws.Range("B3").Value = 'Client-> Create New Client screen : To verify that when the user clicks "Create" button , it navigates to the Create New Client screen'
$ws.Range("C3").Value = "click,click,assert"
$ws.Range("D3").Value = "no value,no value, Create New Client"
$ws.Range("E3").Value = "client_menu,create_contains,client_text"
$ws.Range("F3").Value = "yes"
